$d = $word.ActiveDocument

# 1. Title text change
$d.Content.Find.Execute("Week 2. A Practical Guide to Your Computer", $true, $false, $false, $false, $false, $true, 1, $false, "Week 2. Your Computer", 2)

# 2. Insert a new Heading3 paragraph "Lecture: A Practical Guide to Your Computer"
#    right after the "Readings" heading paragraph, and wrap a new bookmark
#    "X23ecedaa4af68d586a25a58d38f2503d24f4128" starting at that new paragraph
#    and ending at the same place the "readings" bookmark ends (i.e. right
#    after the last bullet of the Readings section).

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Readings", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingPara = $rng.Paragraphs(1)

# Remember where the "readings" bookmark section ends: right after the
# "Note: The Plain Person's..." paragraph, i.e. right before the
# "Data & Computational Work" heading paragraph.
$endRng = $d.Content
$endRng.Find.Execute("Data & Computational Work", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endMarker = $endRng.Paragraphs(1).Previous().Range.End

$headingPara.Range.InsertParagraphAfter()
$newPara = $headingPara.Next()
$newPara.Range.Text = "Lecture: A Practical Guide to Your Computer"
$newPara.Style = "Heading3"

$bmRange = $d.Range($newPara.Range.Start, $endMarker)
$d.Bookmarks.Add("X23ecedaa4af68d586a25a58d38f2503d24f4128", $bmRange)
